# OW-268 fixed notional sign
# Adds the "IRS-Bilateral" worksheet (a bilateral-trade counterpart to
# "IRS-Cleared") and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook
$clearedSheet = $wb.Worksheets.Item("IRS-Cleared")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $clearedSheet)
$ws.Name = "IRS-Bilateral"

# ---- Row 1: headers ----
$ws.Range("A1").Value = 'Value Date'
$ws.Range("B1").Value = 'Position Account ID'
$ws.Range("C1").Value = 'Client ID'
$ws.Range("D1").Value = 'UTI'
$ws.Range("E1").Value = 'Currency'
$ws.Range("F1").Value = 'Effective Date'
$ws.Range("G1").Value = 'Maturity Date'
$ws.Range("H1").Value = 'Cleared Date'
$ws.Range("I1").Value = 'Trade type'
$ws.Range("J1").Value = 'Firm ID'
$ws.Range("K1").Value = 'Source'
$ws.Range("L1").Value = 'LEG1_TYPE'
$ws.Range("M1").Value = 'LEG1_CCY'
$ws.Range("N1").Value = 'LEG1_PAY_FREQ'
$ws.Range("O1").Value = 'LEG1_PAY_ADJ_BUS_DAY_CONV'
$ws.Range("P1").Value = 'LEG1_PAY_ADJ_CAL'
$ws.Range("Q1").Value = 'LEG1_DAYCOUNT'
$ws.Range("R1").Value = 'LEG1_INDEX'
$ws.Range("S1").Value = 'LEG1_INDEX_TENOR'
$ws.Range("T1").Value = 'LEG1_RESET_FREQ'
$ws.Range("U1").Value = 'LEG1_START_DATE'
$ws.Range("V1").Value = 'LEG1_MAT_DATE'
$ws.Range("W1").Value = 'LEG1_NOTIONAL'
$ws.Range("X1").Value = 'LEG1_FIXED_RATE'
$ws.Range("Y1").Value = 'LEG2_TYPE'
$ws.Range("Z1").Value = 'LEG2_CCY'
$ws.Range("AA1").Value = 'LEG2_PAY_FREQ'
$ws.Range("AB1").Value = 'LEG2_PAY_ADJ_BUS_DAY_CONV'
$ws.Range("AC1").Value = 'LEG2_PAY_ADJ_CAL'
$ws.Range("AD1").Value = 'LEG2_DAYCOUNT'
$ws.Range("AE1").Value = 'LEG2_INDEX'
$ws.Range("AF1").Value = 'LEG2_INDEX_TENOR'
$ws.Range("AG1").Value = 'LEG2_RESET_FREQ'
$ws.Range("AH1").Value = 'LEG2_START_DATE'
$ws.Range("AI1").Value = 'LEG2_MAT_DATE'
$ws.Range("AJ1").Value = 'LEG2_NOTIONAL'
$ws.Range("AK1").Value = 'LEG2_FIXED_RATE'
$ws.Range("AL1").Value = 'LEG1_DIRECTION'
$ws.Range("AM1").Value = 'LEG2_DIRECTION'
$ws.Range("AN1").Value = 'Counterpart ID'
$ws.Range("AO1").Value = 'Agreement ID'
$ws.Range("AP1").Value = 'Jurisdiction'

# ---- Row 2: sample bilateral trade data ----
$ws.Range("B2").Value = 'acc1'
$ws.Range("C2").Value = 11811152
$ws.Range("D2").Value = 455820
$ws.Range("E2").Value = 'SGD'
$ws.Range("I2").Value = 'Bilateral'
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 'MARKIT_WIRE'
$ws.Range("L2").Value = 'FIXED'
$ws.Range("M2").Value = 'SGD'
$ws.Range("N2").Value = '6M'
$ws.Range("O2").Value = 'ModifiedFollowing'
$ws.Range("P2").Value = 'SGSI'
$ws.Range("Q2").Value = 'Act/365F'
$ws.Range("W2").Value = '10,000,000.00'
$ws.Range("X2").Value = '1.1'
$ws.Range("Y2").Value = 'FLOAT'
$ws.Range("Z2").Value = 'SGD'
$ws.Range("AA2").Value = '6M'
$ws.Range("AB2").Value = 'ModifiedFollowing'
$ws.Range("AC2").Value = 'SGSI'
$ws.Range("AD2").Value = 'Act/365F'
$ws.Range("AE2").Value = 'SGD-SOR-Reuters'
$ws.Range("AF2").Value = '6M'
$ws.Range("AG2").Value = '6M'
$ws.Range("AJ2").Value = '10,000,000.00'
$ws.Range("AL2").Value = 'R'
$ws.Range("AM2").Value = 'P'
$ws.Range("AN2").Value = 11911171
$ws.Range("AO2").Value = 12011171
$ws.Range("AP2").Value = 'Singapore'

# Date cells use the DD/MM/YY display format (same as on IRS-Cleared)
$ws.Range("A2").NumberFormat = "DD/MM/YY"
$ws.Range("A2").Value = 41631
$ws.Range("F2").NumberFormat = "DD/MM/YY"
$ws.Range("F2").Value = 41607
$ws.Range("G2").NumberFormat = "DD/MM/YY"
$ws.Range("G2").Value = 44164
$ws.Range("H2").NumberFormat = "DD/MM/YY"
$ws.Range("H2").Value = 41605
$ws.Range("U2").NumberFormat = "DD/MM/YY"
$ws.Range("U2").Value = 41607
$ws.Range("V2").NumberFormat = "DD/MM/YY"
$ws.Range("V2").Value = 44164
$ws.Range("AH2").NumberFormat = "DD/MM/YY"
$ws.Range("AH2").Value = 41607
$ws.Range("AI2").NumberFormat = "DD/MM/YY"
$ws.Range("AI2").Value = 44164

# ---- Sheet view / workbook view ----
$ws.Range("A3").Select()

$clearedSheet.Select()
$ws.Activate()
